$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 300
$ws.Range("K2").Value = 300
$ws.Range("M2").Value = -187
$ws.Range("H86").Value = 1010.375
$ws.Range("I86").Value = 964.625
$ws.Range("K86").Value = 964.625
$ws.Range("M86").Value = 158.375
$ws.Range("H89").Value = 1010.375
$ws.Range("I89").Value = 964.625
$ws.Range("K89").Value = 4823.125
$ws.Range("M89").Value = 792.875
$ws.Range("H96").Value = 1628.4286
$ws.Range("I96").Value = 1913
$ws.Range("K96").Value = 5739
$ws.Range("M96").Value = -4366
$ws.Range("H99").Value = 1010.6667
$ws.Range("I99").Value = 434.41666
$ws.Range("J99").Value = 2163.1667
$ws.Range("K99").Value = 1303.24998
$ws.Range("L99").Value = 6489.500100000001
$ws.Range("M99").Value = 194.7500199999999
$ws.Range("N99").Value = -9485.500100000001
$ws.Range("H107").Value = 553.26086
$ws.Range("I107").Value = 371.57895
$ws.Range("J107").Value = 1416.25
$ws.Range("K107").Value = 371.57895
$ws.Range("L107").Value = 1416.25
$ws.Range("M107").Value = 1548.42105
$ws.Range("N107").Value = -5256.25
$ws.Range("H137").Value = 800.56604
$ws.Range("J137").Value = 1409.0714
$ws.Range("L137").Value = 4227.2142
$ws.Range("N137").Value = -9327.2142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 397610.56
$ws.Range("I2").Value = 463751.84
$ws.Range("K2").Value = 463751.84
$ws.Range("M2").Value = -463638.84
$ws.Range("H32").Value = 3132.9058
$ws.Range("I32").Value = 2633.3828
$ws.Range("J32").Value = 13248.25
$ws.Range("K32").Value = 2633.3828
$ws.Range("L32").Value = 13248.25
$ws.Range("M32").Value = -2346.3828
$ws.Range("N32").Value = -13822.25
$ws.Range("H61").Value = 2323.5144
$ws.Range("I61").Value = 1808.2413
$ws.Range("K61").Value = 1808.2413
$ws.Range("M61").Value = -1596.2413
$ws.Range("H97").Value = 869.6818
$ws.Range("I97").Value = 830.9474
$ws.Range("J97").Value = 1115
$ws.Range("K97").Value = 830.9474
$ws.Range("L97").Value = 1115
$ws.Range("M97").Value = -334.9474
$ws.Range("N97").Value = -2107
$ws.Range("H102").Value = 1328.2858
$ws.Range("I102").Value = 1328.2858
$ws.Range("K102").Value = 1328.2858
$ws.Range("M102").Value = 293.7141999999999
$ws.Range("H116").Value = 397610.56
$ws.Range("I116").Value = 463751.84
$ws.Range("K116").Value = 463751.84
$ws.Range("M116").Value = -461457.84
$ws.Range("H132").Value = 1440
$ws.Range("I132").Value = 1117.1086
$ws.Range("J132").Value = 2182.65
$ws.Range("K132").Value = 3351.3258
$ws.Range("L132").Value = 6547.950000000001
$ws.Range("M132").Value = -821.3258000000001
$ws.Range("N132").Value = -11607.95
$ws.Range("H136").Value = 2323.5144
$ws.Range("I136").Value = 1808.2413
$ws.Range("K136").Value = 5424.7239
$ws.Range("M136").Value = -2874.7239

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 903.17645
$ws.Range("I94").Value = 983.73334
$ws.Range("K94").Value = 983.73334
$ws.Range("M94").Value = -532.73334
$ws.Range("H99").Value = 1234.1666
$ws.Range("I99").Value = 1170.2
$ws.Range("K99").Value = 1170.2
$ws.Range("M99").Value = 327.8
$ws.Range("H134").Value = 3786.0476
$ws.Range("I134").Value = 3905.389
$ws.Range("K134").Value = 11716.167
$ws.Range("M134").Value = -9181.167000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1936.5641
$ws.Range("J31").Value = 2562.1177
$ws.Range("L31").Value = 2562.1177
$ws.Range("N31").Value = -3152.1177
$ws.Range("H34").Value = 1936.5641
$ws.Range("J34").Value = 2562.1177
$ws.Range("L34").Value = 2562.1177
$ws.Range("N34").Value = -2966.1177
$ws.Range("H58").Value = 821351.5600000001
$ws.Range("I58").Value = 1175731.8
$ws.Range("J58").Value = 1847.3125
$ws.Range("K58").Value = 1175731.8
$ws.Range("L58").Value = 1847.3125
$ws.Range("M58").Value = -1175528.8
$ws.Range("N58").Value = -2253.3125
$ws.Range("H132").Value = 1263.9714
$ws.Range("I132").Value = 1003.3929
$ws.Range("J132").Value = 2306.2856
$ws.Range("K132").Value = 3010.1787
$ws.Range("L132").Value = 6918.8568
$ws.Range("M132").Value = -480.1787000000004
$ws.Range("N132").Value = -11978.8568
$ws.Range("H134").Value = 1205.7142
$ws.Range("I134").Value = 1089.4193
$ws.Range("J134").Value = 1686.4
$ws.Range("K134").Value = 3268.2579
$ws.Range("L134").Value = 5059.200000000001
$ws.Range("M134").Value = -733.2579000000001
$ws.Range("N134").Value = -10129.2
$ws.Range("H136").Value = 821351.5600000001
$ws.Range("I136").Value = 1175731.8
$ws.Range("J136").Value = 1847.3125
$ws.Range("K136").Value = 3527195.4
$ws.Range("L136").Value = 5541.9375
$ws.Range("M136").Value = -3524645.4
$ws.Range("N136").Value = -10641.9375
$ws.Range("H140").Value = 64500
$ws.Range("J140").Value = 64500
$ws.Range("L140").Value = 64500
$ws.Range("N140").Value = -74860
$ws.Range("H2").Value = 6000
$ws.Range("J2").Value = 6000
$ws.Range("L2").Value = 6000
$ws.Range("N2").Value = -6226
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 85.933334
$ws.Range("I2").Value = 120.55556
$ws.Range("J2").Value = 34
$ws.Range("K2").Value = 723.33336
$ws.Range("L2").Value = 204
$ws.Range("M2").Value = -610.33336
$ws.Range("N2").Value = -430
$ws.Range("H7").Value = 486.76923
$ws.Range("I7").Value = 205
$ws.Range("J7").Value = 728.2857
$ws.Range("K7").Value = 615
$ws.Range("L7").Value = 2184.8571
$ws.Range("M7").Value = -503
$ws.Range("N7").Value = -2408.8571
$ws.Range("H75").Value = 5331.75
$ws.Range("J75").Value = 5331.75
$ws.Range("L75").Value = 15995.25
$ws.Range("N75").Value = -17991.25
$ws.Range("H78").Value = 5331.75
$ws.Range("J78").Value = 5331.75
$ws.Range("L78").Value = 47985.75
$ws.Range("N78").Value = -57969.75
$ws.Range("H118").Value = 55557188
$ws.Range("I118").Value = 333333340
$ws.Range("J118").Value = 1956.6
$ws.Range("K118").Value = 1000000020
$ws.Range("L118").Value = 5869.799999999999
$ws.Range("M118").Value = -999998777
$ws.Range("N118").Value = -8355.799999999999
$ws.Range("H119").Value = 62501376
$ws.Range("I119").Value = 83334500
$ws.Range("J119").Value = 2000
$ws.Range("K119").Value = 250003500
$ws.Range("L119").Value = 6000
$ws.Range("M119").Value = -249998662
$ws.Range("N119").Value = -15676
$ws.Range("H131").Value = 5628012.5
$ws.Range("J131").Value = 10595.619
$ws.Range("L131").Value = 31786.857
$ws.Range("N131").Value = -41866.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 36015.5
$ws.Range("J47").Value = 36015.5
$ws.Range("L47").Value = 36015.5
$ws.Range("N47").Value = -37151.5
$ws.Range("H97").Value = 1115.15
$ws.Range("I97").Value = 947.2353000000001
$ws.Range("J97").Value = 2066.6667
$ws.Range("K97").Value = 947.2353000000001
$ws.Range("L97").Value = 2066.6667
$ws.Range("M97").Value = -451.2353000000001
$ws.Range("N97").Value = -3058.6667
$ws.Range("H132").Value = 803227.2
$ws.Range("I132").Value = 1132889.8
$ws.Range("K132").Value = 3398669.4
$ws.Range("M132").Value = -3396139.4
$ws.Range("H140").Value = 46244.867
$ws.Range("J140").Value = 46244.867
$ws.Range("L140").Value = 46244.867
$ws.Range("N140").Value = -56604.867

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2297.6365
$ws.Range("I22").Value = 2192
$ws.Range("J22").Value = 2482.5
$ws.Range("K22").Value = 2192
$ws.Range("L22").Value = 2482.5
$ws.Range("M22").Value = -1897
$ws.Range("N22").Value = -3072.5
$ws.Range("H27").Value = 2297.6365
$ws.Range("I27").Value = 2192
$ws.Range("J27").Value = 2482.5
$ws.Range("K27").Value = 2192
$ws.Range("L27").Value = 2482.5
$ws.Range("M27").Value = -2085
$ws.Range("N27").Value = -2696.5
$ws.Range("H46").Value = 1714.0834
$ws.Range("I46").Value = 1056
$ws.Range("J46").Value = 2635.4
$ws.Range("K46").Value = 1056
$ws.Range("L46").Value = 2635.4
$ws.Range("M46").Value = -868
$ws.Range("N46").Value = -3011.4
$ws.Range("H93").Value = 2023.25
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 2023.25
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 2023.25
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -4519.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1342.1428
$ws.Range("I81").Value = 732.5
$ws.Range("K81").Value = 1465
$ws.Range("M81").Value = -404
$ws.Range("H84").Value = 1342.1428
$ws.Range("I84").Value = 732.5
$ws.Range("K84").Value = 7325
$ws.Range("M84").Value = -2021
$ws.Range("H96").Value = 10344.333
$ws.Range("I96").Value = 2748.75
$ws.Range("J96").Value = 16420.8
$ws.Range("K96").Value = 2748.75
$ws.Range("L96").Value = 16420.8
$ws.Range("M96").Value = -1375.75
$ws.Range("N96").Value = -19166.8
$ws.Range("H132").Value = 1189.2545
$ws.Range("I132").Value = 871.1707
$ws.Range("J132").Value = 2120.7856
$ws.Range("K132").Value = 2613.5121
$ws.Range("L132").Value = 6362.3568
$ws.Range("M132").Value = -83.51209999999992
$ws.Range("N132").Value = -11422.3568
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
